$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''321.66'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = '''36.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''0.37%'
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''5.135'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''1.92%'
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''0.08111'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''3.61%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''2.161'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''-0.06%'
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''8.011'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''1.34%'
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = 'GateToken'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = '''4.136'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''1.04%'
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = 'MXToken'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = '''0.9287'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''1.29%'
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = '''0.09996'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''2.59%'
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = 'WazirX'
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = '''0.1882'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''1.11%'
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = '''0.09228'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''6.03%'
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = '''0.03592'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''3.33%'
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = '''0.09922'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''0.12%'
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = '''0.001433'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''-0.65%'
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = '''0.005666'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''-0.55%'
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = 'LEO'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = '''3.468'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''0.26%'
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''18.09%'
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = '''-1.50%'
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '''-1.08%'
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = '''6.45%'
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''-0.19%'
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''0.04598'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''-0.21%'
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''0.001244'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''1.30%'
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.004744'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''-7.00%'
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''0.0001301'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''-6.99%'
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''-5.12%'
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = '''0.02031'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''11.25%'
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.05000'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''5.29%'
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.007822'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''2.73%'
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.1401'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''0.46%'
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''0.007818'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''1.01%'
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.002082'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''-6.54%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.01213'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''9.10%'
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.00006423'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''1.03%'
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''0.14%'
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''20.90%'
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''-4.86%'
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.00002101'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''0.14%'
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''0.0002001'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''0.14%'
$ws.Range("E51").Style = "Normal"
